$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.435.59"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.869.42"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.38"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7061"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07902"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3137"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.57"
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07862"
$ws.Range("E11").Value = "  -4.71%  "
$ws.Range("D12").Value = "1.869.94"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.198"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.63"
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7024"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.511"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008365"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "29.434.45"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.84"
$ws.Range("E19").Value = "  +3.52%  "
$ws.Range("D20").Value = "2.137.04"
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.13"
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.639"
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1561"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.010"
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.61"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.83"
$ws.Range("E28").Value = "  +1.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.511"
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.340"
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.266"
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.208"
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05299"
$ws.Range("E33").Value = "  -1.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.899"
$ws.Range("E34").Value = "  -2.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7511"
$ws.Range("E35").Value = "  -1.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.175"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01889"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("D39").Value = "1.281.87"
$ws.Range("E39").Value = "  +2.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.770"
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8951"
$ws.Range("E41").Value = "  -2.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.040"
$ws.Range("E42").Value = "  -7.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "109.02"
$ws.Range("E43").Value = "  -3.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.28"
$ws.Range("E44").Value = "  -3.95%  "
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000128"
$ws.Range("E46").Value = "  -3.57%  "
$ws.Range("D47").Value = "2.035.29"
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.799"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5180"
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.546"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4315"
$ws.Range("E51").Value = "  -1.02%  "